# BDD.xlsx - "Personnes" sheet: add two more player rows (fix for
# "trop de parties affiches si 2 players") so the player-name lookup
# table covers the extra test accounts "wali" and "jimm".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnes")

# Row 11: idPers = 10, nomPers = "wali"
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "wali"
# Keep column C materialised as an (empty) cell, matching the existing rows.
$ws.Cells.Item(11, 3).NumberFormat = "General"

# Row 12: idPers = 11, nomPers = "jimm"
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "jimm"
$ws.Cells.Item(12, 3).NumberFormat = "General"
